# Updated cryptos list on Wed Nov 29 20:50:25 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.772.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "'2.026.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'227.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").Value = "'59.79"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("D10").Value = "'0.0811"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "'14.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "'2.326.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").Value = "'20.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").Value = "'0.756"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "'5.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.32%  "
$ws.Range("D17").Value = "'2.017.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("D18").Value = "'37.666.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").Value = "'6.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("D20").Value = "'69.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "'0.0₃0822"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").Value = "'224.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("D25").Value = "'2.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.69%  "
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").Value = "'165.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("E28").Value = "  -3.82%  "
$ws.Range("D29").Value = "'18.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").Value = "'1.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.59%  "
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("E33").Value = "  +4.38%  "
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("D35").Value = "'0.0601"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.44%  "
$ws.Range("D36").Value = "'6.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.40%  "
$ws.Range("D37").Value = "'2.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.20%  "
$ws.Range("D38").Value = "'3.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.55%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "'1.534.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.45%  "
$ws.Range("D41").Value = "'0.0217"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("D42").Value = "'96.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("D43").Value = "'16.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").Value = "'0.0917"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.33%  "
$ws.Range("D46").Value = "'1.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").Value = "'3.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.53%  "
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("D51").Value = "'2.215.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.58%  "
